$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("streamPWR")

# Select the entire row for "ZERO_CowCr" (row 23) and delete it, shifting
# everything below (including "ZERO_UmpConf") up by one row.
$ws.Range("A23").EntireRow.Select()
$ws.Range("A23").EntireRow.Delete()
